# Integrate esqlabs simulation setup: replace the Aciclovir example scenarios
# on the "Scenarios" sheet with the Amikacin male/female test scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Drop the extra example rows (PopulationScenario, PopulationScenarioFromCSV,
# TestScenario_missingParam) - only two scenario rows remain afterwards.
$ws.Rows("4:6").Delete()

# Clear the two remaining data rows so we can rewrite them from scratch.
$ws.Range("A2:M3").ClearContents()

# New data order matches how the workbook's shared-string table was built
# by the original author (value, then row 3 counterpart, then the names).
$ws.Range("L2").Value = "Amikacin 15mg_kg.pkml"
$ws.Range("B2").Value = "MALE"
$ws.Range("B3").Value = "FEMALE"
$ws.Range("A3").Value = "TestScenarioF"
$ws.Range("A2").Value = "TestScenarioM"
$ws.Range("L3").Value = "Amikacin 15mg_kg.pkml"

# Move the active selection to A2, matching the saved view state.
$ws.Range("A2").Select() | Out-Null
